$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 4 with testmail #2 data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Wil je dit oppakken?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #2: Wil je dit oppakken?"
$logs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F4").Value = "2025-07-29 21:30:57"
$logs.Range("G4").Value = "Nee"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# Extend the conditional formatting ranges from row 2:3 to 2:4
$logs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))
$logs.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H4"))
$logs.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I4"))
$logs.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J4"))

# --- Sheet "Dashboard": append new category row 3 ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B3").Value = 1

# --- Chart: extend category/value series references to include the new row ---
$co = $dash.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
